$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark completed items in the backlog with an "x" in column D
$ws.Range("D2").Value = "x"
$ws.Range("D3").Value = "x"
$ws.Range("D4").Value = "x"
$ws.Range("D9").Value = "x"
$ws.Range("D18").Value = "x"
$ws.Range("D19").Value = "x"

# Update selection to reflect the last edited cell
$ws.Range("D9").Select()
